$d = $word.ActiveDocument

# Original paragraph text is "Version 1." laid out as:
#   run1 = "Version"            (wrapped in a spellStart/spellEnd proofErr pair)
#   run2 = " 1."
#   bookmarkStart/_GoBack/bookmarkEnd
#
# Target paragraph text is "Version 2." laid out as:
#   run1a = "Versi"
#   run1b = "on"
#   run2  = " 2"
#   bookmarkStart/_GoBack/bookmarkEnd
#   run3  = "."

# Step 1: split the "Version" run into "Versi" + "on" (characters 5-7 are "on").
# Re-inserting identical text via InsertXML on that sub-range forces Word to
# break it out of the original run into its own new run, leaving "Versi" behind.
$splitRange = $d.Range(5, 7)
$splitXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>on</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$splitRange.InsertXML($splitXml)

# Step 2: change the version digit "1" to "2" (character 8) in place, which
# keeps it inside the existing " 1." run (now " 2.").
$digitRange = $d.Range(8, 9)
$digitRange.Text = "2"

# Step 3: remove the trailing "." that currently sits before the _GoBack
# bookmark (character 9).
$periodRange = $d.Range(9, 10)
$periodRange.Delete()

# Step 4: re-add the "." as a brand new run positioned after the bookmark
# (i.e. at the current end of the paragraph's text).
$endRange = $d.Range(9, 9)
$endRange.InsertAfter(".")
